$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to Text format so values like "0.4490"
# or "2.200" keep their exact string representation (trailing zeros,
# thousand-dot grouping, etc.) instead of being parsed as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.891.33'
$ws.Range("E2").Value = '  +1.03%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.771.75'
$ws.Range("E3").Value = '  +1.13%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '327.67'
$ws.Range("E5").Value = '  +1.08%  '
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4490'
$ws.Range("E7").Value = '  -2.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3572'
$ws.Range("E8").Value = '  -0.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07468'
$ws.Range("E9").Value = '  -0.60%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.18'
$ws.Range("E10").Value = '  -0.17%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.096'
$ws.Range("E11").Value = '  -0.04%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.003'
$ws.Range("E12").Value = '  +0.20%  '
$ws.Range("E13").Value = '  +0.83%  '
$ws.Range("E14").Value = '  +0.84%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.232'
$ws.Range("E15").Value = '  +1.82%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.770.33'
$ws.Range("E16").Value = '  +1.06%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.99'
$ws.Range("E17").Value = '  +0.62%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001061'
$ws.Range("E18").Value = '  -0.55%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06425'
$ws.Range("E19").Value = '  +0.28%  '
$ws.Range("E20").Value = '  +0.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.23'
$ws.Range("E21").Value = '  +2.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.830'
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.910.61'
$ws.Range("E23").Value = '  +0.92%  '
$ws.Range("E24").Value = '  +1.13%  '
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.88'
$ws.Range("E26").Value = '  -0.75%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.27'
$ws.Range("E27").Value = '  -0.79%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.974.87'
$ws.Range("E28").Value = '  +0.93%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.200'
$ws.Range("E29").Value = '  +5.28%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.89'
$ws.Range("E30").Value = '  -0.65%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.103'
$ws.Range("E31").Value = '  +2.82%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09173'
$ws.Range("E32").Value = '  -0.42%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.585'
$ws.Range("E33").Value = '  +1.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.638'
$ws.Range("E34").Value = '  -0.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.89'
$ws.Range("E35").Value = '  -0.35%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02294'
$ws.Range("E36").Value = '  -0.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06107'
$ws.Range("E37").Value = '  +1.27%  '
$ws.Range("E38").Value = '  -0.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6353'
$ws.Range("E39").Value = '  -0.19%  '
$ws.Range("E40").Value = '  +0.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.187'
$ws.Range("E41").Value = '  -1.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.391'
$ws.Range("E42").Value = '  +0.39%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.925'
$ws.Range("E43").Value = '  +1.75%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.27'
$ws.Range("E44").Value = '  +0.36%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.739'
$ws.Range("E45").Value = '  +0.86%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5882'
$ws.Range("E46").Value = '  -0.38%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '122.63'
$ws.Range("E47").Value = '  -0.31%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.961'
$ws.Range("E48").Value = '  +0.34%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06924'
$ws.Range("E49").Value = '  +1.04%  '
$ws.Range("B50").Value = 'EOS'
$ws.Range("C50").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.141'
$ws.Range("E50").Value = '  -0.41%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.96'
$ws.Range("E51").Value = '  +1.13%  '
